# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" column for the file that was just
# handed off (f7783899-...) on each localized-language status sheet.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet, row 4 = f7783899-... file: new handoff timestamp.
$wsZhCn.Range("D4").Value = "2016-02-24 08:35:54"

# de-de sheet, row 4 = f7783899-... file: new handoff timestamp.
$wsDeDe.Range("D4").Value = "2016-02-24 08:36:04"
